# Applies the "Penalty Reward System" forecast-shift edit:
#   - Forecast Comparison: Week_Start_Date (col B) shifts forward one week
#     for every row, and MyForecast (col D) is overwritten with new values.
#   - Summary: a handful of derived metrics are updated to match.
#
# NOTE: several of these cells hold numbers/dates that are stored as plain
# TEXT in the workbook (t="inlineStr"), not as real numbers/dates. Excel's
# COM layer auto-detects things that *look* like dates/numbers and converts
# them on assignment, so we briefly force the target cell to Text format
# ("@") before writing the literal, then restore "General" so the visible
# formatting is unchanged (matches the original, which has no custom
# number formats).

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison: Week_Start_Date (B) and MyForecast (D) ---
$weekStartDates = @(
    "2025-01-12", "2025-01-19", "2025-01-26", "2025-02-02",
    "2025-02-09", "2025-02-16", "2025-02-23", "2025-03-02",
    "2025-03-09", "2025-03-16", "2025-03-23", "2025-03-30",
    "2025-04-06", "2025-04-13", "2025-04-20", "2025-04-27"
)
$myForecasts = @(4, 4, 3, 6, 7, 3, 3, 3, 3, 3, 3, 3, 3, 3, 3, 3)

for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $i + 2

    $dateCell = $wsForecast.Range("B$row")
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $weekStartDates[$i]
    $dateCell.NumberFormat = "General"

    $wsForecast.Range("D$row").Value = $myForecasts[$i]
}

# --- Summary: update derived metrics (all stored as text) ---
$summaryUpdates = @{
    "B2"  = "2024-02-04 to 2025-01-05";
    "B9"  = "57";
    "B10" = "33";
    "B11" = "18";
    "B12" = "7";
    "B13" = "2025-02-09";
    "B14" = "3";
    "B15" = "2025-01-26";
}

foreach ($addr in $summaryUpdates.Keys) {
    $cell = $wsSummary.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $summaryUpdates[$addr]
    $cell.NumberFormat = "General"
}
